$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column header in H1, copying the formatting (style) used
# by the other header cells (e.g. G1) so it matches the existing header row.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill the new "Save" column (H2:H10) with 0 for every data row.
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}
